# Append a new scraped item "2025-10-14 12:50:45" into the Lancers (ランサーズ)
# sheet's result list. The new item is inserted at row 6 (sorted by score,
# descending) and every existing row from 6 downward is pushed down by one.
# The timestamp column (A) is refreshed to the new run time for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-14 12:50:45"

# --- 1) Shift existing rows 6..15 down to 7..16 (bottom-up so we never
#        clobber a row before it has been read). ---
for ($r = 15; $r -ge 6; $r--) {
    $src = $ws.Range("A" + $r + ":H" + $r).Value()
    $dstRange = $ws.Range("A" + ($r + 1) + ":H" + ($r + 1))
    $dstRange.Value = $src
}

# --- 2) Write the brand-new item into row 6. ---
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【低コスト】住宅リフォーム見積依頼自動化システム構築"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5412955"
$ws.Range("G6").Value = 110
$ws.Range("H6").Value = "◆自動化"

# --- 3) Refresh the "取得日時" timestamp for every data row (2..16). ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("A" + $r).Value = $newTimestamp
}

# --- 4) Rebuild the hyperlinks on column F (2..16) so each URL cell links
#        to the matching address and keeps the Hyperlink style. ---
$ws.Range("F2:F16").Hyperlinks.Delete()
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Range("F" + $r)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
}
